$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-17 (columns A through T), representing the
# 4x4 cross product of sending/target clusters (ECs, FAPs, M2, sCs) for the
# Sema4d-Plxnb1 ligand-receptor pair, per Dr Hou's revised analysis.
$data = @(
    @("ECs", "Sema4d", "Plxnb1", "ECs", 3, 1, 55.59510733333332, 166.785322, 0.4537221086682116, 0.4537221086682116, 2, 0.6666666666666666, 0.7376096666666667, 2.212829, 0.1014284037978316, 0.1014284037978316, 41.00748858843755, 369.067397295938, 0.04602030925000299, 0.04602030925000298),
    @("ECs", "Sema4d", "Plxnb1", "FAPs", 3, 1, 55.59510733333332, 166.785322, 0.4537221086682116, 0.4537221086682116, 3, 1, 0.6291593333333333, 1.887478, 0.08651544278546762, 0.08651544278546762, 34.9781806664351, 314.803625997916, 0.03925396913298639, 0.03925396913298638),
    @("ECs", "Sema4d", "Plxnb1", "M2", 3, 1, 55.59510733333332, 166.785322, 0.4537221086682116, 0.4537221086682116, 1, 0.3333333333333333, 0.07810833333333334, 0.234325, 0.01074064499332162, 0.01074064499332162, 4.342441175294444, 39.08197057765, 0.004873268094826555, 0.004873268094826554),
    @("ECs", "Sema4d", "Plxnb1", "sCs", 3, 1, 55.59510733333332, 166.785322, 0.4537221086682116, 0.4537221086682116, 3, 1, 5.827342666666667, 17.482028, 0.8013155084233792, 0.8013155084233792, 323.9717410214462, 2915.745669193016, 0.3635745621903957, 0.3635745621903956),
    @("FAPs", "Sema4d", "Plxnb1", "ECs", 3, 1, 1.654659, 4.963977, 0.01350398275347337, 0.01350398275347337, 2, 0.6666666666666666, 0.7376096666666667, 2.212829, 0.1014284037978316, 0.1014284037978316, 1.220492473437, 10.984432260933, 0.001369687415598251, 0.001369687415598251),
    @("FAPs", "Sema4d", "Plxnb1", "FAPs", 3, 1, 1.654659, 4.963977, 0.01350398275347337, 0.01350398275347337, 3, 1, 0.6291593333333333, 1.887478, 0.08651544278546762, 0.08651544278546762, 1.041044153334, 9.369397380006, 0.001168303047284067, 0.001168303047284067),
    @("FAPs", "Sema4d", "Plxnb1", "M2", 3, 1, 1.654659, 4.963977, 0.01350398275347337, 0.01350398275347337, 1, 0.3333333333333333, 0.07810833333333334, 0.234325, 0.01074064499332162, 0.01074064499332162, 0.129242656725, 1.163183910525, 0.0001450414847509953, 0.0001450414847509953),
    @("FAPs", "Sema4d", "Plxnb1", "sCs", 3, 1, 1.654659, 4.963977, 0.01350398275347337, 0.01350398275347337, 3, 1, 5.827342666666667, 17.482028, 0.8013155084233792, 0.8013155084233792, 9.642264989484, 86.780384905356, 0.01082095080584006, 0.01082095080584006),
    @("M2", "Sema4d", "Plxnb1", "ECs", 3, 1, 62.05924166666666, 186.177725, 0.5064771225734745, 0.5064771225734744, 2, 0.6666666666666666, 0.7376096666666667, 2.212829, 0.1014284037978316, 0.1014284037978316, 45.7754965593361, 411.979469034025, 0.05137116610274622, 0.0513711661027462),
    @("M2", "Sema4d", "Plxnb1", "FAPs", 3, 1, 62.05924166666666, 186.177725, 0.5064771225734745, 0.5064771225734744, 3, 1, 0.6291593333333333, 1.887478, 0.08651544278546762, 0.08651544278546762, 39.04515111417221, 351.40636002755, 0.04381809252015371, 0.0438180925201537),
    @("M2", "Sema4d", "Plxnb1", "M2", 3, 1, 62.05924166666666, 186.177725, 0.5064771225734745, 0.5064771225734744, 1, 0.3333333333333333, 0.07810833333333334, 0.234325, 0.01074064499332162, 0.01074064499332162, 4.847343934513888, 43.626095410625, 0.005439890970800729, 0.005439890970800728),
    @("M2", "Sema4d", "Plxnb1", "sCs", 3, 1, 62.05924166666666, 186.177725, 0.5064771225734745, 0.5064771225734744, 3, 1, 5.827342666666667, 17.482028, 0.8013155084233792, 0.8013155084233792, 361.6404668251444, 3254.7642014263, 0.4058479729797739, 0.4058479729797738),
    @("sCs", "Sema4d", "Plxnb1", "ECs", 3, 1, 3.222176333333334, 9.666529, 0.02629678600484052, 0.02629678600484052, 2, 0.6666666666666666, 0.7376096666666667, 2.212829, 0.1014284037978316, 0.1014284037978316, 2.376708411171223, 21.390375700541, 0.00266724102948413, 0.00266724102948413),
    @("sCs", "Sema4d", "Plxnb1", "FAPs", 3, 1, 3.222176333333334, 9.666529, 0.02629678600484052, 0.02629678600484052, 3, 1, 0.6291593333333333, 1.887478, 0.08651544278546762, 0.08651544278546762, 2.027262313762444, 18.245360823862, 0.002275078085043466, 0.002275078085043465),
    @("sCs", "Sema4d", "Plxnb1", "M2", 3, 1, 3.222176333333334, 9.666529, 0.02629678600484052, 0.02629678600484052, 1, 0.3333333333333333, 0.07810833333333334, 0.234325, 0.01074064499332162, 0.01074064499332162, 0.2516788231027778, 2.265109407925, 0.0002824444429433403, 0.0002824444429433403),
    @("sCs", "Sema4d", "Plxnb1", "sCs", 3, 1, 3.222176333333334, 9.666529, 0.02629678600484052, 0.02629678600484052, 3, 1, 5.827342666666667, 17.482028, 0.8013155084233792, 0.8013155084233792, 18.77672562675689, 168.990530640812, 0.02107202244736958, 0.02107202244736958)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $targetRow = $startRow + $i
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $ws.Cells.Item($targetRow, $c + 1).Value = $rowValues[$c]
    }
}
